$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy formatting (number formats/styles) from column E into the new column D
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with FY2018 figures (period ending 2018-12-28)
$ws.Range("D7").Value = 43462
$ws.Range("D8").Value = 823600
$ws.Range("D9").Value = 687500
$ws.Range("D10").Value = 136100
$ws.Range("D12").Value = 9400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 15400
$ws.Range("D17").Value = 759600
$ws.Range("D18").Value = 64000
$ws.Range("D20").Value = 200
$ws.Range("D21").Value = 87300
$ws.Range("D22").Value = 10000
$ws.Range("D23").Value = 54200
$ws.Range("D24").Value = -3900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 58100
$ws.Range("D27").Value = 58100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -200
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -200
$ws.Range("D33").Value = 57900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 57900
$ws.Range("D38").Value = 43462
$ws.Range("D41").Value = 43800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 40300
$ws.Range("D44").Value = 121100
$ws.Range("D45").Value = 6300
$ws.Range("D46").Value = 211600
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 41700
$ws.Range("D49").Value = 229900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 485500
$ws.Range("D57").Value = 64300
$ws.Range("D58").Value = 8800
$ws.Range("D59").Value = 14700
$ws.Range("D60").Value = 87800
$ws.Range("D61").Value = 192100
$ws.Range("D62").Value = 7300
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 287200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 59900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 198300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43462
$ws.Range("D81").Value = 57900
$ws.Range("D83").Value = 23100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 60500
$ws.Range("D91").Value = -13900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -15400
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -70600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -25500
